$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.297.06'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = '1.874.30'
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.71'
$ws.Range("E6").Value = '  +1.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3111'
$ws.Range("E8").Value = '  +1.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07744'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08483'
$ws.Range("E11").Value = '  +2.91%  '
$ws.Range("D12").Value = '1.865.81'
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.210'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7113'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.38'
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").Value = '29.298.32'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008307'
$ws.Range("E17").Value = '  +6.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.982'
$ws.Range("E18").Value = '  +2.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.64'
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").Value = '2.126.50'
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.22'
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("E23").Value = '  -1.35%  '
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1624'
$ws.Range("E25").Value = '  +2.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.02'
$ws.Range("E26").Value = '  +0.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.010'
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("E28").Value = '  +2.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.508'
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.420'
$ws.Range("E30").Value = '  +1.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.326'
$ws.Range("E31").Value = '  +6.05%  '
$ws.Range("E32").Value = '  -3.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05259'
$ws.Range("E33").Value = '  +1.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.921'
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7449'
$ws.Range("E36").Value = '  +2.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.684'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01860'
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.718'
$ws.Range("E39").Value = '  +0.84%  '
$ws.Range("D40").Value = '1.162.98'
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.351'
$ws.Range("E41").Value = '  +4.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8898'
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.91'
$ws.Range("E43").Value = '  +1.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '106.86'
$ws.Range("E44").Value = '  +5.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("D46").Value = '2.023.13'
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("E47").Value = '  +2.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000120'
$ws.Range("E49").Value = '  +3.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.370'
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4300'
$ws.Range("E51").Value = '  +1.79%  '
